# Create SheetQueryForm: take the sheet name header block + its query-result
# table and repeat the whole form twice more underneath (one blank separator
# row between each repetition), so the sheet ends up holding the form list
# of queries for the whole sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The header row now also carries a second column value ("insertion") next
# to the sheet/query name.
$ws.Range("B1").Value = "insertion"

# First repetition of the query form block -> rows 6-9 (row 5 left blank).
$ws.Range("A1:B1").Copy($ws.Range("A6"))
$ws.Range("A2:K2").Copy($ws.Range("A7"))
$ws.Range("A3:K3").Copy($ws.Range("A8"))
$ws.Range("A4:K4").Copy($ws.Range("A9"))

# Second repetition of the query form block -> rows 11-14 (row 10 left blank).
$ws.Range("A1:B1").Copy($ws.Range("A11"))
$ws.Range("A2:K2").Copy($ws.Range("A12"))
$ws.Range("A3:K3").Copy($ws.Range("A13"))
$ws.Range("A4:K4").Copy($ws.Range("A14"))

# Leave the selection on the amount cell of the first repeated query.
[void]$ws.Range("K8").Select()
